$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 438.82352
$ws.Range("I28").Value = 492.5
$ws.Range("J28").Value = 310
$ws.Range("K28").Value = 492.5
$ws.Range("L28").Value = 310
$ws.Range("M28").Value = -7.5
$ws.Range("N28").Value = -1280

$ws.Range("H87").Value = 27272
$ws.Range("J87").Value = 27272
$ws.Range("L87").Value = 27272
$ws.Range("N87").Value = -29768

$ws.Range("H90").Value = 27272
$ws.Range("J90").Value = 27272
$ws.Range("L90").Value = 81816
$ws.Range("N90").Value = -94296

$ws.Range("H96").Value = 978
$ws.Range("I96").Value = 745
$ws.Range("J96").Value = 1133.3334
$ws.Range("K96").Value = 2235
$ws.Range("L96").Value = 3400.0002
$ws.Range("M96").Value = -862
$ws.Range("N96").Value = -6146.0002

$ws.Range("H100").Value = 1653.1428
$ws.Range("I100").Value = 1622.7778
$ws.Range("K100").Value = 1622.7778
$ws.Range("M100").Value = -1081.7778

$ws.Range("H137").Value = 3008.1538
$ws.Range("I137").Value = 3536.6
$ws.Range("J137").Value = 2064.5
$ws.Range("K137").Value = 10609.8
$ws.Range("L137").Value = 6193.5
$ws.Range("M137").Value = -8059.799999999999
$ws.Range("N137").Value = -11293.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2731.6667
$ws.Range("I88").Value = 1376
$ws.Range("J88").Value = 3700
$ws.Range("K88").Value = 1376
$ws.Range("L88").Value = 3700
$ws.Range("M88").Value = -970
$ws.Range("N88").Value = -4512

$ws.Range("H91").Value = 2731.6667
$ws.Range("I91").Value = 1376
$ws.Range("J91").Value = 3700
$ws.Range("K91").Value = 1376
$ws.Range("L91").Value = 3700
$ws.Range("M91").Value = 28
$ws.Range("N91").Value = -6508

$ws.Range("H97").Value = 1132.2106
$ws.Range("I97").Value = 348.5
$ws.Range("J97").Value = 3326.6
$ws.Range("K97").Value = 348.5
$ws.Range("L97").Value = 3326.6
$ws.Range("M97").Value = 147.5
$ws.Range("N97").Value = -4318.6

$ws.Range("H110").Value = 1017.0769
$ws.Range("I110").Value = 874.2
$ws.Range("J110").Value = 1493.3334
$ws.Range("K110").Value = 874.2
$ws.Range("L110").Value = 1493.3334
$ws.Range("M110").Value = 1170.8
$ws.Range("N110").Value = -5583.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2963.3845
$ws.Range("I86").Value = 3337.3333
$ws.Range("J86").Value = 2642.8572
$ws.Range("K86").Value = 3337.3333
$ws.Range("L86").Value = 2642.8572
$ws.Range("M86").Value = -2214.3333
$ws.Range("N86").Value = -4888.8572

$ws.Range("H89").Value = 2963.3845
$ws.Range("I89").Value = 3337.3333
$ws.Range("J89").Value = 2642.8572
$ws.Range("K89").Value = 16686.6665
$ws.Range("L89").Value = 13214.286
$ws.Range("M89").Value = -11070.6665
$ws.Range("N89").Value = -24446.286

$ws.Range("H94").Value = 1026.7273
$ws.Range("I94").Value = 1060.9474
$ws.Range("J94").Value = 810
$ws.Range("K94").Value = 1060.9474
$ws.Range("L94").Value = 810
$ws.Range("M94").Value = -609.9474
$ws.Range("N94").Value = -1712

$ws.Range("H134").Value = 4958
$ws.Range("I134").Value = 5544.1665
$ws.Range("J134").Value = 2613.3333
$ws.Range("K134").Value = 16632.4995
$ws.Range("L134").Value = 7839.999899999999
$ws.Range("M134").Value = -14097.4995
$ws.Range("N134").Value = -12909.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 724.5714
$ws.Range("I16").Value = 694.4
$ws.Range("K16").Value = 694.4
$ws.Range("M16").Value = -407.4

$ws.Range("H62").Value = 2600.353
$ws.Range("I62").Value = 2470
$ws.Range("J62").Value = 2786.5715
$ws.Range("K62").Value = 2470
$ws.Range("L62").Value = 2786.5715
$ws.Range("M62").Value = -1846
$ws.Range("N62").Value = -4034.5715

$ws.Range("H65").Value = 2600.353
$ws.Range("I65").Value = 2470
$ws.Range("J65").Value = 2786.5715
$ws.Range("K65").Value = 12350
$ws.Range("L65").Value = 13932.8575
$ws.Range("M65").Value = -9230
$ws.Range("N65").Value = -20172.8575

$ws.Range("H111").Value = 40701.816
$ws.Range("J111").Value = 40701.816
$ws.Range("L111").Value = 40701.816
$ws.Range("N111").Value = -48881.816

$ws.Range("H113").Value = 724.5714
$ws.Range("I113").Value = 694.4
$ws.Range("K113").Value = 694.4
$ws.Range("M113").Value = 1475.6

$ws.Range("H122").Value = 2694.2942
$ws.Range("I122").Value = 2694.2942
$ws.Range("K122").Value = 8082.882599999999
$ws.Range("M122").Value = -5632.882599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5315.346
$ws.Range("I70").Value = 5287
$ws.Range("J70").Value = 5327.9443
$ws.Range("K70").Value = 5287
$ws.Range("L70").Value = 5327.9443
$ws.Range("M70").Value = -5017
$ws.Range("N70").Value = -5867.9443

$ws.Range("H73").Value = 5315.346
$ws.Range("I73").Value = 5287
$ws.Range("J73").Value = 5327.9443
$ws.Range("K73").Value = 5287
$ws.Range("L73").Value = 5327.9443
$ws.Range("M73").Value = -4351
$ws.Range("N73").Value = -7199.9443

$ws.Range("H97").Value = 1745
$ws.Range("I97").Value = 1745
$ws.Range("K97").Value = 1745
$ws.Range("M97").Value = -1249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1465.625
$ws.Range("I61").Value = 1257.7693
$ws.Range("J61").Value = 2366.3333
$ws.Range("K61").Value = 1257.7693
$ws.Range("L61").Value = 2366.3333
$ws.Range("M61").Value = -1055.7693
$ws.Range("N61").Value = -2770.3333

$ws.Range("H113").Value = 1465.625
$ws.Range("I113").Value = 1257.7693
$ws.Range("J113").Value = 2366.3333
$ws.Range("K113").Value = 1257.7693
$ws.Range("L113").Value = 2366.3333
$ws.Range("M113").Value = 912.2307000000001
$ws.Range("N113").Value = -6706.3333

$ws.Range("H141").Value = 76000
$ws.Range("J141").Value = 76000
$ws.Range("L141").Value = 76000
$ws.Range("N141").Value = -86360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6778
$ws.Range("I81").Value = 1217
$ws.Range("J81").Value = 17900
$ws.Range("K81").Value = 2434
$ws.Range("L81").Value = 35800
$ws.Range("M81").Value = -1373
$ws.Range("N81").Value = -37922

$ws.Range("H84").Value = 6778
$ws.Range("I84").Value = 1217
$ws.Range("J84").Value = 17900
$ws.Range("K84").Value = 12170
$ws.Range("L84").Value = 179000
$ws.Range("M84").Value = -6866
$ws.Range("N84").Value = -189608

$ws.Range("H96").Value = 857870.9
$ws.Range("I96").Value = 250774
$ws.Range("J96").Value = 1667333.4
$ws.Range("K96").Value = 250774
$ws.Range("L96").Value = 1667333.4
$ws.Range("M96").Value = -249401
$ws.Range("N96").Value = -1670079.4
